$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Update content on "Admin - Quiz - Edit" (values change from "-" /
#    repeated question text into the new Post-Test placeholder texts)
# ---------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Admin - Quiz - Edit")

$wsEdit.Range("E3").Value = "Post-Test"
$wsEdit.Range("E4").Value = "Belajar apa aja kita?"
$wsEdit.Range("F2").Value = "Lalala"
$wsEdit.Range("G2").Value = "Lalala"
$wsEdit.Range("F3").Value = "Nanana"
$wsEdit.Range("G3").Value = "Nanana"
$wsEdit.Range("F4").Value = "Kakaka"
$wsEdit.Range("G4").Value = "Kakaka"

# ---------------------------------------------------------------------
# 2. Add the new "Admin - Logout" sheet as the last tab
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLogout = $wb.Worksheets.Add($null, $lastSheet)
$wsLogout.Name = "Admin - Logout"

$cCond = $wsLogout.Range("A1")
$cCond.Value = "condition"
$cCond.Font.Bold = $true
$cCond.Borders.ColorIndex = 1
$cCond.Borders.LineStyle = 1
$cCond.HorizontalAlignment = -4108

$cPassed = $wsLogout.Range("A2")
$cPassed.Value = "passed"
$cPassed.Interior.Color = 9359785
$cPassed.Borders.ColorIndex = 1
$cPassed.Borders.LineStyle = 1
$cPassed.HorizontalAlignment = -4108

$wsLogout.Range("B1").Select()

# ---------------------------------------------------------------------
# 3. Move the active tab / selection from "Admin - Quiz - Delete" to
#    "Admin - Quiz - Edit"
# ---------------------------------------------------------------------
$wsEdit.Activate()
$wsEdit.Range("G18").Select()
